# Remove people from the "assignees" sheet who no longer work on the
# Mantid project, then leave the workbook with the "assignees" sheet
# active (mirrors the author switching tabs to review the cleaned-up list).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("assignees")

# Switch to the assignees sheet (this is the sheet the edit was made on).
$ws.Activate()

# Rows (1-based, as originally laid out) of the logins that left the
# project. Delete from the bottom up so earlier indices in this list
# stay valid while we work.
#   2  Pasarus / Sam Jones
#   5  Harrietbrown / Harriet Brown
#   6  martyngigg / Martyn Gigg
#  10  joseph-torsney / Joseph Torsney
#  12  tolu28-coder / Toluwalase Agoro
#  13  DavidFair / David Fair
#  14  StephenSmith25 / Stephen Smith
#  15  DanielMurphy22 / Daniel Murphy
#  17  DannyHindson / Danny Hindson
$rowsToDelete = @(17, 15, 14, 13, 12, 10, 6, 5, 2)
foreach ($r in $rowsToDelete) {
    $ws.Rows.Item($r).Delete()
}

# Leave the selection on the row that was last touched by the delete,
# matching the row where the final removed entry used to sit.
[void]$ws.Rows.Item(9).Select()
